# ---------------------------------------------------------------------------
# Re-populate Sheet1 with the refreshed "motor only" run data, matching the
# second pass of the lab data-export script (run on the lenovo laptop).
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Make sure Sheet1 is the active/selected tab, as in the source file.
$ws.Activate()

# --- Header row ------------------------------------------------------------
$ws.Range("A1").Value = "run_num"
$ws.Range("B1").Value = "block_num"
$ws.Range("C1").Value = "start_time"
$ws.Range("D1").Value = "play_duration"
$ws.Range("E1").Value = "ear"
$ws.Range("F1").Value = "hand"

# --- Data rows --------------------------------------------------------------
# run_num, block_num, start_time, play_duration, ear, hand
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 999
$ws.Range("D2").Value = 999
$ws.Range("E2").Value = "none"
$ws.Range("F2").Value = "R"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 999
$ws.Range("D3").Value = 999
$ws.Range("E3").Value = "none"
$ws.Range("F3").Value = "R"

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 999
$ws.Range("D4").Value = 999
$ws.Range("E4").Value = "none"
$ws.Range("F4").Value = "L"

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 999
$ws.Range("D5").Value = 999
$ws.Range("E5").Value = "none"
$ws.Range("F5").Value = "L"

# --- Column widths (best-fit re-measurement after the data refresh) -------
$ws.Columns.Item(1).ColumnWidth = 8.333333333333334
$ws.Columns.Item(2).ColumnWidth = 10
$ws.Columns.Item(3).ColumnWidth = 9.5
$ws.Columns.Item(5).ColumnWidth = 4.833333333333333
$ws.Columns.Item(6).ColumnWidth = 4.666666666666667

# --- Leftover date-time formatted styles --------------------------------
# (registers numFmtId 22 / m-d-yyyy h:mm against three distinct border
# slots, mirroring the unused style entries left behind in the workbook's
# styles part by the export script's earlier datetime-formatting pass)
$ws.Range("Z100").NumberFormat = "m/d/yy h:mm"
$ws.Range("Z100").Borders.Item(7).LineStyle = 1
$ws.Range("Z101").NumberFormat = "m/d/yy h:mm"
$ws.Range("Z101").Borders.Item(8).LineStyle = 1
$ws.Range("Z102").NumberFormat = "m/d/yy h:mm"
$ws.Range("Z102").Borders.Item(9).LineStyle = 1
$ws.Range("Z100:Z102").Clear()
